{"js": "// Insert a new \"Source: <repo url>\" paragraph right after the\n// \"Author/Engineer: Dennis Komarov\" paragraph and before \"System requirements:\".\n// The new paragraph is bold + single-underlined, matching the other\n// section-header-style paragraphs in this document (e.g. \"System requirements:\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the \"Author/Engineer\" paragraph so the insertion is anchored on\n// content rather than a hard-coded index.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Author/Engineer\") !== -1 && text.indexOf(\"Dennis Komarov\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not locate the \"Author/Engineer: Dennis Komarov\" paragraph.');\n}\n\nconst sourcePara = anchor.insertParagraph(\n  \"Source: https://github.com/311725154/automated_exposure_snap.git\",\n  Word.InsertLocation.after\n);\n\n// Bold + underline the visible run text.\nsourcePara.font.bold = true;\nsourcePara.font.boldBidirectional = true;\nsourcePara.font.underline = Word.UnderlineType.single;\n\n// Also stamp the paragraph mark itself (the end-of-paragraph range) with the\n// same formatting, so the paragraph looks like the other bold/underlined\n// section headers (\"System requirements:\", etc.) even on an empty selection.\nconst endMark = sourcePara.getRange(\"End\");\nendMark.font.bold = true;\nendMark.font.boldBidirectional = true;\nendMark.font.underline = Word.UnderlineType.single;\n\nawait context.sync();\n", "ps1": "# Insert a new \"Source: <repo url>\" paragraph right after the\n# \"Author/Engineer: Dennis Komarov\" paragraph and before \"System requirements:\".\n# The new paragraph is bold + single-underlined, matching the other\n# section-header-style paragraphs in this document (e.g. \"System requirements:\").\n\n$d = $word.ActiveDocument\n\n# Locate the \"Author/Engineer: Dennis Komarov\" paragraph by content (robust to\n# index shifts) instead of hard-coding a paragraph number.\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Author/Engineer*\" -and $t -like \"*Dennis Komarov*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -lt 0) {\n    throw \"Could not locate the 'Author/Engineer: Dennis Komarov' paragraph.\"\n}\n\n# Bail out if the Source paragraph already exists right after it (idempotency).\n$alreadyExists = $false\nif ($targetIndex -lt $count) {\n    $nextText = $d.Paragraphs.Item($targetIndex + 1).Range.Text\n    if ($nextText -like \"*Source:*automated_exposure_snap.git*\") {\n        $alreadyExists = $true\n    }\n}\n\nif (-not $alreadyExists) {\n    $anchorPara = $d.Paragraphs.Item($targetIndex)\n    $anchorPara.Range.InsertParagraphAfter() | Out-Null\n\n    $d2 = $word.ActiveDocument\n    $newPara = $d2.Paragraphs.Item($targetIndex + 1)\n    $r = $newPara.Range\n    $r.Text = \"Source:\"\n    $r.InsertAfter(\" \")\n    $r.InsertAfter(\"https://github.com/311725154/automated_exposure_snap.git\")\n\n    # Bold + single underline, matching the other section-header paragraphs.\n    $r.Font.Bold = 1\n    $r.Font.BoldBi = 1\n    $r.Font.Underline = 1\n}\n"}
